$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the formatting that B4 (hyperlink cell) and E4 (quote-prefixed
# payment-order cell) already carry, so the new/edited cells below can be
# re-stamped with those same style indices instead of whatever new ones
# Range.Value / Hyperlinks.Add would otherwise mint.
$linkStyle = $ws.Range("B4").Style
$numStyle = $ws.Range("E4").Style

# The old hyperlink on B4 is about to move down to B5 - drop it for now,
# we'll re-add hyperlinks (to B5 and B4) further below.
$ws.Range("B4").Hyperlinks.Delete()

# --- New row 5: a copy of what used to be in row 4 (ocerutti / preprod),
#     just with an updated payment-order number.
$ws.Range("A5").Value = $ws.Range("A4").Value()
$ws.Range("B5").Value = $ws.Range("B4").Value()
$ws.Range("C5").Value = "ocerutti"
$ws.Range("D5").Value = $ws.Range("D4").Value()
$ws.Range("E5").Style = $numStyle
$ws.Range("E5").Value = "'1120170200936   "

# --- Row 4: now holds the dgariffo user with the new payment-order number.
$ws.Range("C4").Value = "dgariffo"
$ws.Range("E4").Value = "'1220170301429   "

# --- B5 keeps the hyperlink that used to live on B4; B4 gets a freshly
#     entered hyperlink to the same URL.
$ws.Hyperlinks.Add($ws.Range("B5"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do") | Out-Null
$ws.Range("B4").Style = $linkStyle
$ws.Range("B5").Style = $linkStyle

# --- Row 2/3: E2's order number changed; E3 keeps the same text (its
#     shared-string index simply shifts because of the table edit above).
$ws.Range("E2").Value = "'1220170301429"
$ws.Range("E3").Value = "'1120194100412"

# --- Misc view state captured in the diff.
$ws.Range("D8").Select() | Out-Null
